# Auto-generated edit script applying the committed value changes
# to the Excalibur_Profits workbook (sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 977.25
$ws.Range("J17").Value = 984.6829
$ws.Range("L17").Value = 2954.0487
$ws.Range("N17").Value = -3290.0487
$ws.Range("H33").Value = 416.42856
$ws.Range("I33").Value = 240.8125
$ws.Range("K33").Value = 240.8125
$ws.Range("M33").Value = -11.8125
$ws.Range("H86").Value = 2208.0667
$ws.Range("I86").Value = 1902
$ws.Range("K86").Value = 1902
$ws.Range("M86").Value = -779
$ws.Range("H89").Value = 2208.0667
$ws.Range("I89").Value = 1902
$ws.Range("K89").Value = 9510
$ws.Range("M89").Value = -3894
$ws.Range("H101").Value = 1245.9231
$ws.Range("I101").Value = 981.5454999999999
$ws.Range("J101").Value = 2700
$ws.Range("K101").Value = 2944.6365
$ws.Range("L101").Value = 8100
$ws.Range("M101").Value = -1322.6365
$ws.Range("N101").Value = -11344
$ws.Range("H106").Value = 2943.2354
$ws.Range("I106").Value = 2824.5715
$ws.Range("K106").Value = 2824.5715
$ws.Range("M106").Value = -2193.5715
$ws.Range("H112").Value = 9267.182000000001
$ws.Range("I112").Value = 349.5
$ws.Range("J112").Value = 9842.517
$ws.Range("K112").Value = 1048.5
$ws.Range("L112").Value = 29527.551
$ws.Range("M112").Value = 59.5
$ws.Range("N112").Value = -31743.551
$ws.Range("H116").Value = 248450.36
$ws.Range("I116").Value = 98870.5
$ws.Range("J116").Value = 330039.38
$ws.Range("K116").Value = 98870.5
$ws.Range("L116").Value = 330039.38
$ws.Range("M116").Value = -95428.5
$ws.Range("N116").Value = -336923.38
$ws.Range("H132").Value = 198894
$ws.Range("I132").Value = 215157.19
$ws.Range("K132").Value = 645471.5700000001
$ws.Range("M132").Value = -642941.5700000001
$ws.Range("H137").Value = 816838.9
$ws.Range("J137").Value = 1052214.8
$ws.Range("L137").Value = 3156644.4
$ws.Range("N137").Value = -3161744.4
$ws.Range("H138").Value = 1971.2941
$ws.Range("I138").Value = 1275.2858
$ws.Range("J138").Value = 3095.6155
$ws.Range("K138").Value = 3825.8574
$ws.Range("L138").Value = 9286.8465
$ws.Range("M138").Value = 1314.1426
$ws.Range("N138").Value = -19566.8465

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2298
$ws.Range("I45").Value = 2348.5
$ws.Range("J45").Value = 2197
$ws.Range("K45").Value = 2348.5
$ws.Range("L45").Value = 2197
$ws.Range("M45").Value = -1971.5
$ws.Range("N45").Value = -2951
$ws.Range("H61").Value = 2231209.5
$ws.Range("I61").Value = 2390083.8
$ws.Range("K61").Value = 2390083.8
$ws.Range("M61").Value = -2389871.8
$ws.Range("H74").Value = 8266.454
$ws.Range("I74").Value = 10262.071
$ws.Range("J74").Value = 4774.125
$ws.Range("K74").Value = 10262.071
$ws.Range("L74").Value = 4774.125
$ws.Range("M74").Value = -9388.071
$ws.Range("N74").Value = -6522.125
$ws.Range("H77").Value = 8266.454
$ws.Range("I77").Value = 10262.071
$ws.Range("J77").Value = 4774.125
$ws.Range("K77").Value = 51310.355
$ws.Range("L77").Value = 23870.625
$ws.Range("M77").Value = -46942.355
$ws.Range("N77").Value = -32606.625
$ws.Range("H132").Value = 1978229.8
$ws.Range("I132").Value = 1978229.8
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 5934689.4
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -5932159.4
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 2231209.5
$ws.Range("I136").Value = 2390083.8
$ws.Range("K136").Value = 7170251.399999999
$ws.Range("M136").Value = -7167701.399999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3166.4443
$ws.Range("I99").Value = 4166.3335
$ws.Range("J99").Value = 1166.6666
$ws.Range("K99").Value = 4166.3335
$ws.Range("L99").Value = 1166.6666
$ws.Range("M99").Value = -2668.3335
$ws.Range("N99").Value = -4162.6666
$ws.Range("H107").Value = 6001.8335
$ws.Range("I107").Value = 6001.8335
$ws.Range("K107").Value = 6001.8335
$ws.Range("M107").Value = -4081.8335
$ws.Range("H134").Value = 3552497.5
$ws.Range("I134").Value = 23809524
$ws.Range("J134").Value = 658636.5600000001
$ws.Range("K134").Value = 71428572
$ws.Range("L134").Value = 1975909.68
$ws.Range("M134").Value = -71426037
$ws.Range("N134").Value = -1980979.68

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 6328.3335
$ws.Range("I62").Value = 1992.5
$ws.Range("J62").Value = 15000
$ws.Range("K62").Value = 1992.5
$ws.Range("L62").Value = 15000
$ws.Range("M62").Value = -1368.5
$ws.Range("N62").Value = -16248
$ws.Range("H65").Value = 6328.3335
$ws.Range("I65").Value = 1992.5
$ws.Range("J65").Value = 15000
$ws.Range("K65").Value = 9962.5
$ws.Range("L65").Value = 75000
$ws.Range("M65").Value = -6842.5
$ws.Range("N65").Value = -81240
$ws.Range("H99").Value = 18262.834
$ws.Range("I99").Value = 1889
$ws.Range("K99").Value = 1889
$ws.Range("M99").Value = -391
$ws.Range("H126").Value = 18262.834
$ws.Range("I126").Value = 1889
$ws.Range("K126").Value = 5667
$ws.Range("M126").Value = -3197
$ws.Range("H132").Value = 50264984
$ws.Range("I132").Value = 90921190
$ws.Range("K132").Value = 272763570
$ws.Range("M132").Value = -272761040

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 8189.7144
$ws.Range("I3").Value = 7465.6
$ws.Range("K3").Value = 22396.8
$ws.Range("M3").Value = -22284.8
$ws.Range("H23").Value = 71.82353000000001
$ws.Range("I23").Value = 48.5
$ws.Range("K23").Value = 145.5
$ws.Range("M23").Value = 89.5
$ws.Range("H34").Value = 1398.7778
$ws.Range("J34").Value = 2874.75
$ws.Range("L34").Value = 8624.25
$ws.Range("N34").Value = -8792.25
$ws.Range("H114").Value = 2262.15
$ws.Range("J114").Value = 2464.7222
$ws.Range("L114").Value = 7394.1666
$ws.Range("N114").Value = -13902.1666
$ws.Range("H137").Value = 1925.3334
$ws.Range("J137").Value = 3350
$ws.Range("L137").Value = 10050
$ws.Range("N137").Value = -20250

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 6422
$ws.Range("I31").Value = 6422
$ws.Range("K31").Value = 6422
$ws.Range("M31").Value = -6130
$ws.Range("H37").Value = 6422
$ws.Range("I37").Value = 6422
$ws.Range("K37").Value = 6422
$ws.Range("M37").Value = -6145
$ws.Range("H80").Value = 372066.16
$ws.Range("J80").Value = 4761
$ws.Range("L80").Value = 4761
$ws.Range("N80").Value = -6757
$ws.Range("H83").Value = 372066.16
$ws.Range("J83").Value = 4761
$ws.Range("L83").Value = 23805
$ws.Range("N83").Value = -33789
$ws.Range("H102").Value = 5790.1055
$ws.Range("I102").Value = 5633.7666
$ws.Range("K102").Value = 5633.7666
$ws.Range("M102").Value = -4011.7666
$ws.Range("H126").Value = 1114140.1
$ws.Range("I126").Value = 1854311.5
$ws.Range("K126").Value = 5562934.5
$ws.Range("M126").Value = -5560464.5
$ws.Range("H132").Value = 18748980
$ws.Range("I132").Value = 28117798
$ws.Range("K132").Value = 84353394
$ws.Range("M132").Value = -84350864

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3033.3333
$ws.Range("I61").Value = 2050
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 2050
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1848
$ws.Range("N61").Value = -5404
$ws.Range("H100").Value = 16743.428
$ws.Range("I100").Value = 2866.6667
$ws.Range("K100").Value = 2866.6667
$ws.Range("M100").Value = -2325.6667
$ws.Range("H113").Value = 3033.3333
$ws.Range("I113").Value = 2050
$ws.Range("J113").Value = 5000
$ws.Range("K113").Value = 2050
$ws.Range("L113").Value = 5000
$ws.Range("M113").Value = 120
$ws.Range("N113").Value = -9340
$ws.Range("H132").Value = 1586399
$ws.Range("I132").Value = 2322718.8
$ws.Range("K132").Value = 6968156.399999999
$ws.Range("M132").Value = -6965626.399999999
$ws.Range("H136").Value = 63305.715
$ws.Range("I136").Value = 2999.8572
$ws.Range("J136").Value = 93458.64
$ws.Range("K136").Value = 8999.571599999999
$ws.Range("L136").Value = 280375.92
$ws.Range("M136").Value = -6449.571599999999
$ws.Range("N136").Value = -285475.92

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1104.1177
$ws.Range("I100").Value = 784.7692
$ws.Range("K100").Value = 1569.5384
$ws.Range("M100").Value = -1028.5384
$ws.Range("H107").Value = 1854.2667
$ws.Range("I107").Value = 1363.3158
$ws.Range("K107").Value = 4089.9474
$ws.Range("M107").Value = -2169.9474
$ws.Range("H113").Value = 6390.875
$ws.Range("J113").Value = 7796.4
$ws.Range("L113").Value = 23389.2
$ws.Range("N113").Value = -27729.2
$ws.Range("H132").Value = 3149271.5
$ws.Range("I132").Value = 4110693.2
$ws.Range("J132").Value = 8627.333000000001
$ws.Range("K132").Value = 12332079.6
$ws.Range("L132").Value = 25881.999
$ws.Range("M132").Value = -12329549.6
$ws.Range("N132").Value = -30941.999

